# Revert capacity chart to show kilowatts on the y-axis:
#  - Divide the previously-misscaled Watt values (Solar, Wind series) by 1000
#    so they read as kilowatts.
#  - Update the numeric display format for the data cells to show one
#    decimal place (kW values are no longer whole numbers).
#  - Update the chart's value-axis title and number format accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Correct the capacity values (Watts -> Kilowatts) -------------------
# Solar column (E): rows 17, 21-26
$ws.Range("E17").Value = 8
$ws.Range("E21").Value = 13.4
$ws.Range("E22").Value = 11
$ws.Range("E23").Value = 13.2
$ws.Range("E24").Value = 27.9
$ws.Range("E25").Value = 27.88
$ws.Range("E26").Value = 25.06

# Wind column (G): row 11
$ws.Range("G11").Value = 35.4

# --- 2. Data cells now carry one decimal place ------------------------------
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3. Update the chart's value axis title + number format ----------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

$valAxis = $chart.Axes(2)  # xlValue
$valAxis.AxisTitle.Text = "Kilowatts (kW)"
$valAxis.TickLabels.NumberFormat = "#,##0"
